$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'246.55"
$ws.Cells.Item(3, 4).Value = "'29.70"
$ws.Cells.Item(3, 5).Value = "'9.25%"
$ws.Cells.Item(4, 4).Value = "'5.173"
$ws.Cells.Item(4, 5).Value = "'1.30%"
$ws.Cells.Item(5, 4).Value = "'0.05731"
$ws.Cells.Item(5, 5).Value = "'0.43%"
$ws.Cells.Item(6, 5).Value = "'1.16%"
$ws.Cells.Item(7, 4).Value = "'0.8556"
$ws.Cells.Item(7, 5).Value = "'4.40%"
$ws.Cells.Item(8, 4).Value = "'0.8674"
$ws.Cells.Item(8, 5).Value = "'1.15%"
$ws.Cells.Item(9, 4).Value = "'0.1365"
$ws.Cells.Item(9, 5).Value = "'2.67%"
$ws.Cells.Item(10, 4).Value = "'0.07075"
$ws.Cells.Item(10, 5).Value = "'1.76%"
$ws.Cells.Item(11, 4).Value = "'0.02929"
$ws.Cells.Item(11, 5).Value = "'3.06%"
$ws.Cells.Item(12, 4).Value = "'0.09381"
$ws.Cells.Item(12, 5).Value = "'-0.20%"
$ws.Cells.Item(13, 4).Value = "'0.001512"
$ws.Cells.Item(13, 5).Value = "'-0.38%"
$ws.Cells.Item(14, 4).Value = "'0.04141"
$ws.Cells.Item(14, 5).Value = "'2.93%"
$ws.Cells.Item(15, 4).Value = "'0.0006006"
$ws.Cells.Item(15, 5).Value = "'0.48%"
$ws.Cells.Item(16, 4).Value = "'0.006178"
$ws.Cells.Item(16, 5).Value = "'0.21%"
$ws.Cells.Item(17, 5).Value = "'5,070.11%"
$ws.Cells.Item(18, 4).Value = "'3.489"
$ws.Cells.Item(18, 5).Value = "'-0.65%"
$ws.Cells.Item(19, 5).Value = "'2.83%"
$ws.Cells.Item(20, 4).Value = "'2.274"
$ws.Cells.Item(20, 5).Value = "'-1.86%"
$ws.Cells.Item(21, 4).Value = "'0.3155"
$ws.Cells.Item(21, 5).Value = "'-0.29%"
$ws.Cells.Item(22, 4).Value = "'0.03442"
$ws.Cells.Item(22, 5).Value = "'6.90%"
$ws.Cells.Item(23, 4).Value = "'0.1303"
$ws.Cells.Item(23, 5).Value = "'0.00%"
$ws.Cells.Item(24, 4).Value = "'3.463"
$ws.Cells.Item(24, 5).Value = "'-2.87%"
$ws.Cells.Item(25, 4).Value = "'0.1379"
$ws.Cells.Item(25, 5).Value = "'0.41%"
$ws.Cells.Item(27, 5).Value = "'0.27%"
$ws.Cells.Item(40, 4).Value = "'0.03753"
$ws.Cells.Item(40, 5).Value = "'0.71%"
$ws.Cells.Item(41, 2).Value = "KickToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Cells.Item(41, 4).Value = "'0.005736"
$ws.Cells.Item(41, 5).Value = "'-3.39%"
$ws.Cells.Item(42, 2).Value = "BKEXToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Cells.Item(42, 4).Value = "'0.1073"
$ws.Cells.Item(42, 5).Value = "'1.39%"
$ws.Cells.Item(43, 2).Value = "CEJI"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Cells.Item(43, 4).Value = "'0.002299"
$ws.Cells.Item(43, 5).Value = "'-0.04%"
$ws.Cells.Item(44, 4).Value = "'0.008496"
$ws.Cells.Item(44, 5).Value = "'-10.94%"
$ws.Cells.Item(45, 4).Value = "'0.00005258"
$ws.Cells.Item(45, 5).Value = "'2.34%"
$ws.Cells.Item(46, 5).Value = "'-0.04%"
$ws.Cells.Item(47, 4).Value = "'0.06466"
$ws.Cells.Item(47, 5).Value = "'-35.96%"
$ws.Cells.Item(48, 4).Value = "'0.002519"
$ws.Cells.Item(48, 5).Value = "'0.05%"
$ws.Cells.Item(49, 5).Value = "'-0.04%"
$ws.Cells.Item(50, 5).Value = "'-0.04%"
